# Updates the cryptocurrency price/volume table on the active sheet to
# reflect the refreshed scrape (GitHub Actions "Updated cryptos list" run).
#
# Most rows keep their Coin/Link in place and only get new Price (D) and
# Volume(1h) (E) values. A handful of rows only changed one of the two
# columns. Rows 50/51 (Cronos / Mantle) swapped rank order, so those two
# rows get their Coin/Link/Price/Volume all rewritten together.
#
# All the source values (prices like "1.750.85", "0.000008589", etc.) are
# plain text in the workbook, not real numbers - Coinranking formats them
# with dots as thousand separators and Excel would otherwise happily
# "helpfully" reinterpret numeric-looking strings (losing precision, or
# flipping to scientific notation for tiny values). Force every written
# cell to stay text, matching the existing (unstyled) text cells.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both Price (D) and Volume(1h) (E) changed.
$updates = @(
    @{ Row = 2; D = "27.944.36"; E = "  -4.13%  " },
    @{ Row = 3; D = "1.741.76"; E = "  -4.46%  " },
    @{ Row = 4; D = "1.007"; E = "  +0.35%  " },
    @{ Row = 5; D = "225.54"; E = "  -3.55%  " },
    @{ Row = 6; D = "0.5776"; E = "  -3.23%  " },
    @{ Row = 7; D = "1.006"; E = "  +0.20%  " },
    @{ Row = 8; D = "0.2697"; E = "  -1.70%  " },
    @{ Row = 9; D = "22.87"; E = "  -0.98%  " },
    @{ Row = 10; D = "0.06564"; E = "  -5.39%  " },
    @{ Row = 11; D = "0.07501"; E = "  -1.25%  " },
    @{ Row = 12; D = "1.757.14"; E = "  -3.83%  " },
    @{ Row = 13; D = "4.685"; E = "  -0.90%  " },
    @{ Row = 14; D = "0.6004"; E = "  -3.46%  " },
    @{ Row = 15; D = "1.980.54"; E = "  -4.49%  " },
    @{ Row = 16; D = "73.59"; E = "  -5.77%  " },
    @{ Row = 17; D = "0.000008577"; E = "  -10.70%  " },
    @{ Row = 18; D = "28.000.15"; E = "  -2.15%  " },
    @{ Row = 19; D = "5.267"; E = "  -5.43%  " },
    @{ Row = 20; D = "1.007"; E = "  +0.31%  " },
    @{ Row = 21; D = "204.26"; E = "  -5.68%  " },
    @{ Row = 22; D = "11.21"; E = "  -2.49%  " },
    @{ Row = 23; D = "6.615"; E = "  -2.81%  " },
    @{ Row = 24; D = "1.008"; E = "  +0.28%  " },
    @{ Row = 25; D = "150.43"; E = "  -4.02%  " },
    @{ Row = 26; D = "7.965"; E = "  +0.61%  " },
    @{ Row = 27; D = "0.1225"; E = "  -4.28%  " },
    @{ Row = 28; D = "15.97"; E = "  -2.87%  " },
    @{ Row = 29; D = "1.405"; E = "  -0.68%  " },
    @{ Row = 30; D = "0.06015"; E = "  -4.25%  " },
    @{ Row = 31; D = "1.382"; E = "  -3.73%  " },
    @{ Row = 32; D = "3.718"; E = "  -2.58%  " },
    @{ Row = 33; D = "3.688"; E = "  -1.31%  " },
    @{ Row = 34; D = "1.658"; E = "  -3.31%  " },
    @{ Row = 35; D = "1.025"; E = "  -5.27%  " },
    @{ Row = 36; D = "0.6271"; E = "  -2.35%  " },
    @{ Row = 37; D = "2.372"; E = "  -6.46%  " },
    @{ Row = 39; D = "6.236"; E = "  -5.04%  " },
    @{ Row = 40; D = "0.01650"; E = "  -5.66%  " },
    @{ Row = 41; D = "1.123.33"; E = "  -1.72%  " },
    @{ Row = 42; D = "0.8569"; E = "  -3.15%  " },
    @{ Row = 44; D = "99.27"; E = "  -0.99%  " },
    @{ Row = 46; D = "58.63"; E = "  -5.00%  " },
    @{ Row = 48; D = "1.558"; E = "  -1.95%  " },
    @{ Row = 49; D = "8.250"; E = "  -2.30%  " }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
}

# Rows where only the Volume(1h) (E) column changed.
Set-TextValue $ws.Cells.Item(38, 5) "  -3.90%  "
Set-TextValue $ws.Cells.Item(43, 5) "  +0.49%  "
Set-TextValue $ws.Cells.Item(47, 5) "  -3.81%  "

# Row where only the Price (D) column changed.
Set-TextValue $ws.Cells.Item(45, 4) "1.894.38"

# Rows 50/51: Cronos and Mantle swapped positions in the ranking, each
# with its own refreshed Price/Volume values.
Set-TextValue $ws.Cells.Item(50, 2) "Mantle"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Cells.Item(50, 4) "0.4455"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.99%  "

Set-TextValue $ws.Cells.Item(51, 2) "Cronos"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(51, 4) "0.05392"
Set-TextValue $ws.Cells.Item(51, 5) "  -2.16%  "
